$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("M4").Value = 1.11
$ws.Range("N4").Value = 6.5

$ws.Range("M5").Value = 1.05
$ws.Range("N5").Value = 11

$ws.Range("G6").Value = 2.7
$ws.Range("I6").Value = 2.6
$ws.Range("J6").Value = 3.25
$ws.Range("L6").Value = 3.2
$ws.Range("X6").Value = 13
$ws.Range("Z6").Value = 26
$ws.Range("AJ6").Value = 26
$ws.Range("AK6").Value = 21
$ws.Range("AW6").Value = 4.75
$ws.Range("AX6").Value = 15

$ws.Range("Q8").Value = 2.15
$ws.Range("R8").Value = 1.67

$ws.Range("U9").Value = 1.91
$ws.Range("V9").Value = 1.8

$ws.Range("G10").Value = 2.38
$ws.Range("H10").Value = 2.75
$ws.Range("I10").Value = 3.6
$ws.Range("L10").Value = 4
$ws.Range("O10").Value = 1.53
$ws.Range("P10").Value = 2.38
$ws.Range("Q10").Value = 2.7
$ws.Range("R10").Value = 1.44
$ws.Range("V10").Value = 1.62
$ws.Range("AG10").Value = 8
$ws.Range("AI10").Value = 13
$ws.Range("AR10").Value = 81

$ws.Range("U11").Value = 1.73

$ws.Range("U12").Value = 2.63
$ws.Range("V12").Value = 1.44

$wb.Save()
